# Auto-generated Excel COM-interop script to apply Titan_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 39750
$ws.Range("J13").Value = 36333.332
$ws.Range("L13").Value = 36333.332
$ws.Range("N13").Value = -36671.332
$ws.Range("H107").Value = 695013.5600000001
$ws.Range("I107").Value = 855267.25
$ws.Range("K107").Value = 855267.25
$ws.Range("M107").Value = -853347.25
$ws.Range("H113").Value = 1956.1666
$ws.Range("I113").Value = 1967.4
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1967.4
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 1286.6
$ws.Range("N113").Value = -8408
$ws.Range("H132").Value = 267048.97
$ws.Range("I132").Value = 312212.25
$ws.Range("J132").Value = 15424.857
$ws.Range("K132").Value = 936636.75
$ws.Range("L132").Value = 46274.571
$ws.Range("M132").Value = -934106.75
$ws.Range("N132").Value = -51334.571
$ws.Range("H133").Value = 41109.75
$ws.Range("J133").Value = 41109.75
$ws.Range("L133").Value = 41109.75
$ws.Range("N133").Value = -51229.75
$ws.Range("H134").Value = 71393.336
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 71393.336
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 71393.336
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -81533.336
$ws.Range("H136").Value = 58945
$ws.Range("J136").Value = 58945
$ws.Range("L136").Value = 58945
$ws.Range("N136").Value = -69145
$ws.Range("H137").Value = 34484228
$ws.Range("I137").Value = 55556584
$ws.Range("J137").Value = 2193.4546
$ws.Range("K137").Value = 166669752
$ws.Range("L137").Value = 6580.3638
$ws.Range("M137").Value = -166667202
$ws.Range("N137").Value = -11680.3638
$ws.Range("H138").Value = 8611840
$ws.Range("I138").Value = 2528425
$ws.Range("J138").Value = 12823435
$ws.Range("K138").Value = 7585275
$ws.Range("L138").Value = 38470305
$ws.Range("M138").Value = -7580135
$ws.Range("N138").Value = -38480585
$ws.Range("H140").Value = 53536
$ws.Range("J140").Value = 53536
$ws.Range("L140").Value = 53536
$ws.Range("N140").Value = -63896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3548.4866
$ws.Range("I132").Value = 3055
$ws.Range("J132").Value = 5083.778
$ws.Range("K132").Value = 9165
$ws.Range("L132").Value = 15251.334
$ws.Range("M132").Value = -6635
$ws.Range("N132").Value = -20311.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 40000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 40000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 40000
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -40504
$ws.Range("H107").Value = 1224.6875
$ws.Range("I107").Value = 903
$ws.Range("J107").Value = 1760.8334
$ws.Range("K107").Value = 903
$ws.Range("L107").Value = 1760.8334
$ws.Range("M107").Value = 1017
$ws.Range("N107").Value = -5600.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2718.4119
$ws.Range("I31").Value = 1718.5
$ws.Range("J31").Value = 4146.857
$ws.Range("K31").Value = 1718.5
$ws.Range("L31").Value = 4146.857
$ws.Range("M31").Value = -1423.5
$ws.Range("N31").Value = -4736.857
$ws.Range("H34").Value = 2718.4119
$ws.Range("I34").Value = 1718.5
$ws.Range("J34").Value = 4146.857
$ws.Range("K34").Value = 1718.5
$ws.Range("L34").Value = 4146.857
$ws.Range("M34").Value = -1516.5
$ws.Range("N34").Value = -4550.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 378.42105
$ws.Range("I6").Value = 114.46154
$ws.Range("J6").Value = 950.3333
$ws.Range("K6").Value = 343.38462
$ws.Range("L6").Value = 2850.9999
$ws.Range("M6").Value = -230.38462
$ws.Range("N6").Value = -3076.9999
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 300
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -638
$ws.Range("H26").Value = 394.5
$ws.Range("I26").Value = 216.66667
$ws.Range("J26").Value = 501.2
$ws.Range("K26").Value = 650.00001
$ws.Range("L26").Value = 1503.6
$ws.Range("M26").Value = -362.00001
$ws.Range("N26").Value = -2079.6
$ws.Range("H33").Value = 52.214287
$ws.Range("I33").Value = 53.625
$ws.Range("K33").Value = 321.75
$ws.Range("M33").Value = -38.75
$ws.Range("H41").Value = 845.7143
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 784
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 2352
$ws.Range("M41").Value = -2662
$ws.Range("N41").Value = -3028
$ws.Range("H63").Value = 4870.6665
$ws.Range("I63").Value = 4870.6665
$ws.Range("K63").Value = 14611.9995
$ws.Range("M63").Value = -13862.9995
$ws.Range("H66").Value = 4870.6665
$ws.Range("I66").Value = 4870.6665
$ws.Range("K66").Value = 43835.9985
$ws.Range("M66").Value = -40091.9985
$ws.Range("H119").Value = 481.6
$ws.Range("I119").Value = 481.6
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 1444.8
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 3393.2
$ws.Range("N119").ClearContents()
$ws.Range("H131").Value = 1306.0154
$ws.Range("I131").Value = 573.3333
$ws.Range("J131").Value = 1380.5254
$ws.Range("K131").Value = 1719.9999
$ws.Range("L131").Value = 4141.5762
$ws.Range("M131").Value = 3320.0001
$ws.Range("N131").Value = -14221.5762

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5917.143
$ws.Range("I70").Value = 5730.909
$ws.Range("J70").Value = 6600
$ws.Range("K70").Value = 5730.909
$ws.Range("L70").Value = 6600
$ws.Range("M70").Value = -5460.909
$ws.Range("N70").Value = -7140
$ws.Range("H73").Value = 5917.143
$ws.Range("I73").Value = 5730.909
$ws.Range("J73").Value = 6600
$ws.Range("K73").Value = 5730.909
$ws.Range("L73").Value = 6600
$ws.Range("M73").Value = -4794.909
$ws.Range("N73").Value = -8472
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H138").Value = 67241.63
$ws.Range("J138").Value = 67241.63
$ws.Range("L138").Value = 67241.63
$ws.Range("N138").Value = -77521.63
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 966.6667
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 950
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 950
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1540
$ws.Range("H27").Value = 966.6667
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 950
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 950
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1164
$ws.Range("H35").Value = 15614.286
$ws.Range("I35").Value = 20660
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 20660
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = -20324
$ws.Range("N35").Value = -3672
$ws.Range("H39").Value = 40032.5
$ws.Range("I39").Value = 50000
$ws.Range("J39").Value = 30065
$ws.Range("K39").Value = 50000
$ws.Range("L39").Value = 30065
$ws.Range("M39").Value = -49540
$ws.Range("N39").Value = -30985
$ws.Range("H58").Value = 6500
$ws.Range("I58").Value = 6500
$ws.Range("K58").Value = 6500
$ws.Range("M58").Value = -6240
$ws.Range("H132").Value = 3870.3914
$ws.Range("I132").Value = 2312.7856
$ws.Range("J132").Value = 6293.3335
$ws.Range("K132").Value = 6938.3568
$ws.Range("L132").Value = 18880.0005
$ws.Range("M132").Value = -4408.3568
$ws.Range("N132").Value = -23940.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 6166
$ws.Range("I17").Value = 7749
$ws.Range("J17").Value = 3000
$ws.Range("K17").Value = 7749
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -7577
$ws.Range("N17").Value = -3344
$ws.Range("H122").Value = 84845.586
$ws.Range("I122").Value = 112239.11
$ws.Range("J122").Value = 2665
$ws.Range("K122").Value = 336717.33
$ws.Range("L122").Value = 7995
$ws.Range("M122").Value = -334267.33
$ws.Range("N122").Value = -12895
